# Weekly update: a new price record (for 2022-09-15 / serial 44819) is
# reported for "Comercializadora del Agro de Limarí - Arándano (blue)".
# It becomes the new first data row; every existing record shifts down
# one row (dimension grows from A1:T11 to A1:T12).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 2, pushing all data rows
# (old rows 2-11) down to rows 3-12.
$ws.Rows("2:2").Insert()

# The inserted row inherits the header row's bold/centered formatting;
# reset it, then restore the date-time number format on D2 to match the
# other "Fecha" cells in column D.
$ws.Rows("2:2").ClearFormats()
$ws.Range("D2").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Populate the new row with the latest weekly record.
$ws.Range("A2").Value = 2
$ws.Range("B2").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C2").Value = "Coquimbo"
$ws.Range("D2").Value = 44819
$ws.Range("E2").Value = 4
$ws.Range("F2").Value = "Fruta"
$ws.Range("G2").Value = 100101
$ws.Range("H2").Value = "Berries"
$ws.Range("I2").Value = 100101001
$ws.Range("J2").Value = "Arándano (blue)"
$ws.Range("K2").Value = "Sin especificar"
$ws.Range("L2").Value = "Primera"
$ws.Range("M2").Value = 240
$ws.Range("N2").Value = 11000
$ws.Range("O2").Value = 12000
$ws.Range("P2").Value = 11500
$ws.Range("Q2").Value = "$/bandeja 2 kilos"
$ws.Range("R2").Value = "Provincia de Limarí"
$ws.Range("S2").Value = 5750
$ws.Range("T2").Value = 2
